$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.259.85"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "3.057.37"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.33"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.77"
$ws.Range("E6").Value = "  +5.69%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.547"
$ws.Range("E8").Value = "  +3.29%  "
$ws.Range("D9").Value = "3.069.82"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.85"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.60"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "3.563.55"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "63.295.44"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "3.065.72"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.29"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.67"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.99"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.04"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  +3.82%  "
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.25"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "0.0₃0850"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.11"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  +3.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.23"
$ws.Range("E39").Value = "  -3.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.33"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.57"
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "444.50"
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.286"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0363"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("E45").Value = "  +3.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.10"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").Value = "2.802.58"
$ws.Range("E47").Value = "  -4.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.28"
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.36"
$ws.Range("E49").Value = "  +4.84%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.27"
$ws.Range("E51").Value = "  +0.87%  "
